$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => (Price, Volume(1h)) updates. $null means leave that column unchanged.
# Price values that look like plain numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the source data, which
# keeps Price as a text column, e.g. "30.633.10" / "249.89").
$updates = @{
    2  = @("30.633.10", "  +0.30%  ")
    3  = @("1.962.29", "  +2.27%  ")
    4  = @($null, "  -0.16%  ")
    5  = @("'249.89", "  +1.75%  ")
    6  = @("'0.9993", "  -0.18%  ")
    7  = @("'0.4835", "  +0.88%  ")
    8  = @("'0.2946", "  +2.04%  ")
    9  = @("'0.06796", "  +1.11%  ")
    10 = @("'110.42", "  +0.69%  ")
    11 = @($null, "  +1.28%  ")
    12 = @("1.987.06", "  +3.70%  ")
    13 = @("'0.07746", "  +2.31%  ")
    14 = @("'5.465", "  +3.79%  ")
    15 = @("'0.6887", "  +3.12%  ")
    16 = @("'295.37", "  -0.78%  ")
    17 = @("30.657.85", "  +0.49%  ")
    18 = @("2.247.05", "  +4.08%  ")
    19 = @($null, "  +1.80%  ")
    20 = @("'0.000007711", "  +1.88%  ")
    21 = @("'5.620", "  +0.35%  ")
    22 = @("'0.9993", "  -0.15%  ")
    23 = @("'0.9996", "  -0.21%  ")
    24 = @("'6.623", "  +3.10%  ")
    25 = @("'9.926", "  +4.64%  ")
    26 = @("'170.71", "  +3.68%  ")
    27 = @("'20.14", "  -0.67%  ")
    28 = @("'2.222", "  +5.11%  ")
    29 = @("'0.1064", "  -1.23%  ")
    30 = @("'1.436", "  +2.90%  ")
    31 = @("'4.727", "  +17.22%  ")
    32 = @($null, "  +6.88%  ")
    33 = @("'0.05128", "  +2.63%  ")
    34 = @("'0.7736", "  +5.01%  ")
    35 = @("'1.179", "  +3.60%  ")
    36 = @("'0.02050", "  +0.64%  ")
    37 = @("'2.729", $null)
    38 = @("'2.716", "  +1.15%  ")
    39 = @("'2.109", "  +4.39%  ")
    40 = @("'6.343", "  +7.28%  ")
    41 = @("'0.4481", "  +1.17%  ")
    42 = @("'109.61", "  -1.03%  ")
    43 = @("'0.8762", "  +1.77%  ")
    44 = @("'70.42", "  -3.38%  ")
    45 = @("'0.9999", "  -0.08%  ")
    46 = @("'7.474", "  +2.82%  ")
    47 = @("'0.1283", "  +4.17%  ")
    48 = @("'9.397", "  +0.95%  ")
    49 = @("'36.12", "  +2.93%  ")
    50 = @("'47.73", "  -3.17%  ")
    51 = @("'0.4091", "  +2.19%  ")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($null -ne $price) {
        $ws.Cells.Item($row, 4).Value = $price
    }
    if ($null -ne $volume) {
        $ws.Cells.Item($row, 5).Value = $volume
    }
}
